$wb = $excel.ActiveWorkbook

# The workbook keeps duplicate "想去人数" (F column) figures on both the
# 展览 sheet and the combined 全部类型 sheet. Bump the figures on both.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F13").Value = 517
    $ws.Range("F15").Value = 13787
    $ws.Range("F18").Value = 9078
    $ws.Range("F20").Value = 8192
}
